$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# Sheet "Weekly Quantity": append rows 28-30
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(28, 45662.99999999999, 2),
    @(29, 45669.99999999999, 2),
    @(30, 45683.99999999999, 1)
)

foreach ($row in $weeklyNewRows) {
    $r = $row[0]
    $wsWeekly.Range("A$r").NumberFormat = $dateFormat
    $wsWeekly.Range("A$r").Value = $row[1]
    $wsWeekly.Range("B$r").Value = $row[2]
}

# ---------------------------------------------------------------------------
# Sheet "Monthly Trend": append row 16
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Range("A16").NumberFormat = $dateFormat
$wsMonthly.Range("A16").Value = 45688.99999999999
$wsMonthly.Range("B16").Value = 5

# ---------------------------------------------------------------------------
# Sheet "PO Forecast": update existing rows 2-35 and append rows 36-38
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Rows whose B value changes but A (date) stays the same
$forecastBUpdates = @(
    @(2, 3),
    @(3, 3),
    @(4, 3),
    @(5, 3),
    @(8, 4),
    @(9, 4),
    @(10, 4),
    @(18, 7),
    @(19, 7),
    @(20, 7),
    @(21, 7),
    @(22, 7),
    @(23, 7),
    @(24, 7),
    @(25, 8),
    @(26, 8),
    @(27, 8)
)

foreach ($row in $forecastBUpdates) {
    $r = $row[0]
    $wsForecast.Range("B$r").Value = $row[1]
}

# Rows 28-38: both the date (A) and the value (B) are rewritten - the
# forecast series shifts forward and extends with three new future rows.
$forecastFull = @(
    @(28, 45662.99999999999, 9),
    @(29, 45669.99999999999, 9),
    @(30, 45683.99999999999, 9),
    @(31, 45690.99999999999, 9),
    @(32, 45697.99999999999, 9),
    @(33, 45704.99999999999, 9),
    @(34, 45711.99999999999, 9),
    @(35, 45718.99999999999, 9),
    @(36, 45725.99999999999, 9),
    @(37, 45732.99999999999, 9),
    @(38, 45739.99999999999, 9)
)

foreach ($row in $forecastFull) {
    $r = $row[0]
    $wsForecast.Range("A$r").NumberFormat = $dateFormat
    $wsForecast.Range("A$r").Value = $row[1]
    $wsForecast.Range("B$r").Value = $row[2]
}
